$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2..B5 content (role/content table rows), and add a new row 6.
# Values are assigned in the order their text first appears in the shared
# string table of the target workbook, so the new <si> entries line up.
$ws.Range("B2").Value = "Eu sou um assistente de ajuda chamado Fuleco em homenagem ao mascote da copa do mundo do Brasil."

$ws.Range("A6").Value = "system"
$ws.Range("B6").Value = "Eu sou brasileiro e falo apenas português ao responder alguém"

$ws.Range("B4").Value = "Uso primeira pessoa ao escrever minhas frases"

$ws.Range("B3").Value = "Sempre falo do meu nome ao começar uma frase e sempre em tom muito alegre"

$ws.Range("B5").Value = "Eu explico minuciosamente todas as perguntas"

# Rows 1-5 use a custom height of 18.75, the new row 6 keeps the default height.
$ws.Rows.Item(4).RowHeight = 18.75
$ws.Rows.Item(5).RowHeight = 18.75

# Widen column B to fit the longer text (approx. 123.4 chars) and drop best-fit.
$ws.Columns.Item(2).ColumnWidth = 122.6

# Leave B3 selected, matching the end state captured in the saved file.
$ws.Range("B3").Select() | Out-Null
